$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the shared strings in the order they first appear in the target file so
# the resulting sharedStrings table matches: 0=peroos, 1=INICIO, 2=FIN
$ws.Range("A2").Value = "peroos"
$ws.Range("A1").Value = "INICIO"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "peroos"
    $ws.Cells.Item($r, 2).Value = 65799
    $ws.Cells.Item($r, 3).Value = 123456
}

$ws.Range("A12").Value = "FIN"
$ws.Cells.Item(12, 2).Value = 65799
$ws.Cells.Item(12, 3).Value = 123456

# Underline style applied to a couple of otherwise empty cells
$ws.Range("F9").Font.Underline = $true
$ws.Range("E16").Font.Underline = $true

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("E14").Select() | Out-Null
